$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of K2:L16 (text_bdrc_id / text_84000_ids data rows)
$ws.Range("K2:L16").ClearContents()

# Select the range that was cleared, matching the resulting selection state
$ws.Range("K2:L16").Select()
